$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite columns B:E (rows 1-3) with values copied from columns O, R, AN, AQ
# (per commit "Hjemme passive tweaks lichtwark deleted values")
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 5.9891685282022342
$ws.Range("C2").Value = 4.740854810894489
$ws.Range("D2").Value = 7.7660648385154882
$ws.Range("E2").Value = 7.1047710993981532

$ws.Range("B3").Value = 4.8852490163363234
$ws.Range("C3").Value = 7.8360232974744992
$ws.Range("D3").Value = 8.6431981979258197
$ws.Range("E3").Value = 6.4173446571310562

# Update the selected range to reflect the new, smaller region of interest
$ws.Range("B1:E3").Select()
